# Auto-generated edit script applying the Ragnarok_Profits market-data refresh
# (per-row currentAveragePrice / Leve profit recalculation) described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 21962.531
$ws.Range("J17").Value = 21962.531
$ws.Range("L17").Value = 65887.59299999999
$ws.Range("N17").Value = -66223.59299999999

$ws.Range("H38").Value = 2559.5833
$ws.Range("I38").Value = 314.75
$ws.Range("J38").Value = 7049.25
$ws.Range("K38").Value = 944.25
$ws.Range("L38").Value = 21147.75
$ws.Range("M38").Value = -572.25
$ws.Range("N38").Value = -21891.75

$ws.Range("H41").Value = 1289.875
$ws.Range("J41").Value = 267
$ws.Range("L41").Value = 267
$ws.Range("N41").Value = -1147

$ws.Range("H62").Value = 5248.5454
$ws.Range("I62").Value = 2441.75
$ws.Range("K62").Value = 2441.75
$ws.Range("M62").Value = -1817.75

$ws.Range("H65").Value = 5248.5454
$ws.Range("I65").Value = 2441.75
$ws.Range("K65").Value = 12208.75
$ws.Range("M65").Value = -9088.75

$ws.Range("H97").Value = 1990.125
$ws.Range("J97").Value = 1990.125
$ws.Range("L97").Value = 5970.375
$ws.Range("N97").Value = -6962.375

$ws.Range("H98").Value = 3572763.8
$ws.Range("I98").Value = 4630911
$ws.Range("K98").Value = 4630911
$ws.Range("M98").Value = -4629413

$ws.Range("H100").Value = 8698.223
$ws.Range("I100").Value = 5875.8
$ws.Range("J100").Value = 12226.25
$ws.Range("K100").Value = 5875.8
$ws.Range("L100").Value = 12226.25
$ws.Range("M100").Value = -5334.8
$ws.Range("N100").Value = -13308.25

$ws.Range("H101").Value = 2786.3845
$ws.Range("J101").Value = 5199.8335
$ws.Range("L101").Value = 15599.5005
$ws.Range("N101").Value = -18843.5005

$ws.Range("H106").Value = 9242.0625
$ws.Range("I106").Value = 8682.538
$ws.Range("K106").Value = 8682.538
$ws.Range("M106").Value = -8051.538

$ws.Range("H122").Value = 3572763.8
$ws.Range("I122").Value = 4630911
$ws.Range("K122").Value = 13892733
$ws.Range("M122").Value = -13890283

$ws.Range("H135").Value = 1619.5883
$ws.Range("I135").Value = 494.5357
$ws.Range("K135").Value = 4450.821300000001
$ws.Range("M135").Value = -1915.821300000001

$ws.Range("H137").Value = 684264.3
$ws.Range("I137").Value = 667.1818
$ws.Range("J137").Value = 1624210.4
$ws.Range("K137").Value = 2001.5454
$ws.Range("L137").Value = 4872631.199999999
$ws.Range("M137").Value = 548.4546
$ws.Range("N137").Value = -4877731.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3401.9517
$ws.Range("I32").Value = 3406.2222
$ws.Range("K32").Value = 3406.2222
$ws.Range("M32").Value = -3119.2222

$ws.Range("H61").Value = 14010091
$ws.Range("I61").Value = 20006246
$ws.Range("J61").Value = 2017779.8
$ws.Range("K61").Value = 20006246
$ws.Range("L61").Value = 2017779.8
$ws.Range("M61").Value = -20006034
$ws.Range("N61").Value = -2018203.8

$ws.Range("H74").Value = 834902.5
$ws.Range("I74").Value = 927224.2
$ws.Range("K74").Value = 927224.2
$ws.Range("M74").Value = -926350.2

$ws.Range("H77").Value = 834902.5
$ws.Range("I77").Value = 927224.2
$ws.Range("K77").Value = 4636121
$ws.Range("M77").Value = -4631753

$ws.Range("H92").Value = 87226.5
$ws.Range("J92").Value = 87226.5
$ws.Range("L92").Value = 87226.5
$ws.Range("N92").Value = -92218.5

$ws.Range("H110").Value = 4951.6665
$ws.Range("I110").Value = 3942
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 3942
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = -1897
$ws.Range("N110").Value = -14090

$ws.Range("H132").Value = 2945685
$ws.Range("I132").Value = 4495.8076
$ws.Range("J132").Value = 12504550
$ws.Range("K132").Value = 13487.4228
$ws.Range("L132").Value = 37513650
$ws.Range("M132").Value = -10957.4228
$ws.Range("N132").Value = -37518710

$ws.Range("H136").Value = 14010091
$ws.Range("I136").Value = 20006246
$ws.Range("J136").Value = 2017779.8
$ws.Range("K136").Value = 60018738
$ws.Range("L136").Value = 6053339.4
$ws.Range("M136").Value = -60016188
$ws.Range("N136").Value = -6058439.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1093.1
$ws.Range("I22").Value = 1364
$ws.Range("K22").Value = 1364
$ws.Range("M22").Value = -1191

$ws.Range("H86").Value = 346522.66
$ws.Range("I86").Value = 509227.75
$ws.Range("J86").Value = 6321.091
$ws.Range("K86").Value = 509227.75
$ws.Range("L86").Value = 6321.091
$ws.Range("M86").Value = -508104.75
$ws.Range("N86").Value = -8567.091

$ws.Range("H89").Value = 346522.66
$ws.Range("I89").Value = 509227.75
$ws.Range("J89").Value = 6321.091
$ws.Range("K89").Value = 2546138.75
$ws.Range("L89").Value = 31605.455
$ws.Range("M89").Value = -2540522.75
$ws.Range("N89").Value = -42837.455

$ws.Range("H105").Value = 2080204.8
$ws.Range("I105").Value = 2855981.5
$ws.Range("K105").Value = 2855981.5
$ws.Range("M105").Value = -2854234.5

$ws.Range("H107").Value = 2449
$ws.Range("I107").Value = 1768.6666
$ws.Range("J107").Value = 2959.25
$ws.Range("K107").Value = 1768.6666
$ws.Range("L107").Value = 2959.25
$ws.Range("M107").Value = 151.3334
$ws.Range("N107").Value = -6799.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33672624
$ws.Range("I31").Value = 52634120
$ws.Range("J31").Value = 920943.75
$ws.Range("K31").Value = 52634120
$ws.Range("L31").Value = 920943.75
$ws.Range("M31").Value = -52633825
$ws.Range("N31").Value = -921533.75

$ws.Range("H34").Value = 33672624
$ws.Range("I34").Value = 52634120
$ws.Range("J34").Value = 920943.75
$ws.Range("K34").Value = 52634120
$ws.Range("L34").Value = 920943.75
$ws.Range("M34").Value = -52633918
$ws.Range("N34").Value = -921347.75

$ws.Range("H62").Value = 41673024
$ws.Range("I62").Value = 3813
$ws.Range("K62").Value = 3813
$ws.Range("M62").Value = -3189

$ws.Range("H65").Value = 41673024
$ws.Range("I65").Value = 3813
$ws.Range("K65").Value = 19065
$ws.Range("M65").Value = -15945

$ws.Range("H105").Value = 1736.3334
$ws.Range("I105").Value = 1736.3334
$ws.Range("K105").Value = 1736.3334
$ws.Range("M105").Value = 10.66660000000002

$ws.Range("H107").Value = 5599.5
$ws.Range("I107").Value = 1199.5
$ws.Range("J107").Value = 9999.5
$ws.Range("K107").Value = 1199.5
$ws.Range("L107").Value = 9999.5
$ws.Range("M107").Value = 720.5
$ws.Range("N107").Value = -13839.5

$ws.Range("H132").Value = 2840.625
$ws.Range("I132").Value = 2810.3157
$ws.Range("K132").Value = 8430.947100000001
$ws.Range("M132").Value = -5900.947100000001

$ws.Range("H134").Value = 3284.0454
$ws.Range("I134").Value = 3047.889
$ws.Range("K134").Value = 9143.667000000001
$ws.Range("M134").Value = -6608.667000000001

$ws.Range("H141").Value = 172468.44
$ws.Range("J141").Value = 172468.44
$ws.Range("L141").Value = 172468.44
$ws.Range("N141").Value = -182828.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1469.7333
$ws.Range("J12").Value = 2157.2
$ws.Range("L12").Value = 6471.599999999999
$ws.Range("N12").Value = -6817.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6585683
$ws.Range("I132").Value = 2776.0454
$ws.Range("J132").Value = 27274820
$ws.Range("K132").Value = 8328.136200000001
$ws.Range("L132").Value = 81824460
$ws.Range("M132").Value = -5798.136200000001
$ws.Range("N132").Value = -81829520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 110001
$ws.Range("I22").Value = 110001
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 110001
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -109706

$ws.Range("N27").ClearContents()
$ws.Range("H27").Value = 110001
$ws.Range("I27").Value = 110001
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 110001
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -109894

$ws.Range("H68").Value = 5210868
$ws.Range("I68").Value = 10418541
$ws.Range("K68").Value = 10418541
$ws.Range("M68").Value = -10417792

$ws.Range("H71").Value = 5210868
$ws.Range("I71").Value = 10418541
$ws.Range("K71").Value = 52092705
$ws.Range("M71").Value = -52088961

$ws.Range("H82").Value = 9500
$ws.Range("J82").Value = 18000
$ws.Range("L82").Value = 18000
$ws.Range("N82").Value = -18722

$ws.Range("H85").Value = 9500
$ws.Range("J85").Value = 18000
$ws.Range("L85").Value = 18000
$ws.Range("N85").Value = -20496

$ws.Range("H132").Value = 4869.5
$ws.Range("I132").Value = 3022.625
$ws.Range("K132").Value = 9067.875
$ws.Range("M132").Value = -6537.875

$ws.Range("H136").Value = 4038.5293
$ws.Range("I136").Value = 3040.9092
$ws.Range("K136").Value = 9122.7276
$ws.Range("M136").Value = -6572.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M31").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0

$ws.Range("H56").Value = 82745
$ws.Range("I56").Value = 63285
$ws.Range("J56").Value = 92475
$ws.Range("K56").Value = 63285
$ws.Range("L56").Value = 92475
$ws.Range("M56").Value = -62571
$ws.Range("N56").Value = -93903

$ws.Range("H126").Value = 8571.615
$ws.Range("I126").Value = 9339.046
$ws.Range("J126").Value = 4350.75
$ws.Range("K126").Value = 28017.138
$ws.Range("L126").Value = 13052.25
$ws.Range("M126").Value = -25547.138
$ws.Range("N126").Value = -17992.25

$ws.Range("H132").Value = 457055.8
$ws.Range("I132").Value = 2671.1875
$ws.Range("K132").Value = 8013.5625
$ws.Range("M132").Value = -5483.5625
